$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.912.01"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.649.10"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.80%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.87"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3838"
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.26"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.343"
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08445"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.81"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.012"
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.925"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001314"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").Value = "1.649.59"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.95"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06963"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.50"
$ws.Range("E20").Value = "  -3.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.941"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.61"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").Value = "23.865.91"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.443"
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.908"
$ws.Range("E26").Value = "  -4.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.93"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.01"
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.384"
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "137.34"
$ws.Range("E30").Value = "  -2.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.709"
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.487"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").Value = "1.830.45"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08113"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9876"
$ws.Range("E35").Value = "  -5.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.699"
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02918"
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2680"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.50"
$ws.Range("E39").Value = "  -3.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09113"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7544"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.43"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.421"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.73"
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6920"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.438"
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.090"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08266"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.42"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("E51").Value = "  -1.54%  "
